$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-08-07 Monday" "2023-08-08 Tuesday"

Replace-Text "31×38=1178" "90×95=8550"
Replace-Text "90×74=6660" "24×22=528"
Replace-Text "36×18=648" "72×94=6768"
Replace-Text "49×38=1862" "58×14=812"
Replace-Text "47×56=2632" "62×48=2976"
Replace-Text "44×68=2992" "13×20=260"
Replace-Text "15×25=375" "77×71=5467"
Replace-Text "42×37=1554" "91×45=4095"
Replace-Text "55×37=2035" "97×58=5626"
Replace-Text "25×68=1700" "92×99=9108"
Replace-Text "29×48=1392" "51×86=4386"
Replace-Text "98×13=1274" "29×59=1711"
Replace-Text "60×94=5640" "14×63=882"
Replace-Text "60×18=1080" "89×58=5162"
Replace-Text "61×17=1037" "32×21=672"
Replace-Text "85×52=4420" "96×49=4704"
Replace-Text "98×89=8722" "62×57=3534"
Replace-Text "16×45=720" "42×31=1302"
Replace-Text "39×50=1950" "19×29=551"
Replace-Text "90×83=7470" "49×81=3969"
Replace-Text "40×59=2360" "25×65=1625"
Replace-Text "68×32=2176" "47×12=564"
Replace-Text "85×67=5695" "43×32=1376"
Replace-Text "47×74=3478" "11×60=660"
Replace-Text "38×67=2546" "79×36=2844"

Write-Output "Done"
